$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @"
questions = [
    {
        "title": "I put people under pressure.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I joke around a lot.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I hesitate to criticize other people's ideas.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I am emotionally reserved.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    }
]
"@

# Put the new pretty-printed JSON-style content into A2 (it already
# carries no special formatting), then delete row 1 (the old A1=0
# header cell with the bold/border style) so A2 shifts up to A1.
$ws.Range("A2").Value = $text
$ws.Rows.Item(1).Delete() | Out-Null

# Writing a multi-line value auto-expands the row height; re-fit the
# row so it collapses back to the sheet's default height.
$ws.Rows.Item(1).EntireRow.AutoFit() | Out-Null
